$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.497.01'
$ws.Range("E2").Value = '  -0.84%  '

$ws.Range("D3").Value = '1.849.87'
$ws.Range("E3").Value = '  -0.38%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9993'
$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.09'
$ws.Range("E5").Value = '  -0.62%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6363'
$ws.Range("E6").Value = '  -1.08%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '47.65'
$ws.Range("E8").Value = '  +0.43%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07557'
$ws.Range("E9").Value = '  +0.78%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.2988'
$ws.Range("E10").Value = '  +0.13%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '24.19'
$ws.Range("E11").Value = '  -1.26%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07686'
$ws.Range("E12").Value = '  +0.38%  '

$ws.Range("D13").Value = '1.875.41'
$ws.Range("E13").Value = '  +0.72%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.024'

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6866'
$ws.Range("E15").Value = '  -0.55%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '83.90'
$ws.Range("E16").Value = '  -0.11%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000009731'
$ws.Range("E17").Value = '  +1.13%  '

$ws.Range("D18").Value = '2.115.16'
$ws.Range("E18").Value = '  +0.22%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.224'
$ws.Range("E19").Value = '  +2.60%  '

$ws.Range("D20").Value = '29.548.47'
$ws.Range("E20").Value = '  -0.67%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '236.53'
$ws.Range("E21").Value = '  +0.34%  '

$ws.Range("E22").Value = '  -1.06%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.0000'
$ws.Range("E23").Value = '  +0.02%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.628'
$ws.Range("E24").Value = '  +2.16%  '

$ws.Range("E25").Value = '  +0.00%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '155.92'
$ws.Range("E26").Value = '  -1.73%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1391'
$ws.Range("E27").Value = '  -2.16%  '

$ws.Range("E28").Value = '  -0.80%  '

$ws.Range("E29").Value = '  -0.90%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.488'
$ws.Range("E30").Value = '  -0.45%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.05864'
$ws.Range("E31").Value = '  -6.92%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.272'
$ws.Range("E32").Value = '  -0.63%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.124'
$ws.Range("E33").Value = '  -0.67%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.058'
$ws.Range("E34").Value = '  -0.79%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.901'
$ws.Range("E35").Value = '  +0.18%  '

$ws.Range("E36").Value = '  -0.29%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.7174'
$ws.Range("E37").Value = '  -1.68%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.596'
$ws.Range("E38").Value = '  -0.40%  '

$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.800'
$ws.Range("E39").Value = '  -1.53%  '

$ws.Range("B40").Value = 'Maker'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D40").Value = '1.237.09'
$ws.Range("E40").Value = '  +3.00%  '

$ws.Range("E41").Value = '  -0.80%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9146'
$ws.Range("E42").Value = '  -0.73%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.116'
$ws.Range("E43").Value = '  -0.70%  '

$ws.Range("E44").Value = '  -0.05%  '

$ws.Range("D45").Value = '2.030.23'
$ws.Range("E45").Value = '  +0.44%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '67.50'

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '101.76'
$ws.Range("E47").Value = '  -0.34%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.346'
$ws.Range("E48").Value = '  +9.78%  '

$ws.Range("B49").Value = 'TheSandbox'
$ws.Range("C49").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4037'
$ws.Range("E49").Value = '  -0.76%  '

$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.146'
$ws.Range("E50").Value = '  -0.92%  '

$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.697'
$ws.Range("E51").Value = '  +2.64%  '
